# Auto-generated edit script: updates cryptos list values per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# are pre-formatted as Text so the value round-trips as a string, matching
# the source data (which stores these as plain text, e.g. "44.101.14").
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "44.101.14"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "2.242.08"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "315.80"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "99.04"
$ws.Range("E6").Value = "  -5.97%  "
$ws.Range("D7").Value = "0.576"
$ws.Range("E7").Value = "  -2.92%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  -6.58%  "
$ws.Range("D10").Value = "36.06"
$ws.Range("E10").Value = "  -6.75%  "
$ws.Range("E11").Value = "  -2.48%  "
$ws.Range("E12").Value = "  -6.65%  "
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("D14").Value = "2.581.96"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").Value = "0.844"
$ws.Range("E15").Value = "  -4.53%  "
$ws.Range("D16").Value = "2.242.60"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").Value = "13.92"
$ws.Range("E17").Value = "  -4.37%  "
$ws.Range("D18").Value = "43.913.29"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").Value = "13.11"
$ws.Range("E19").Value = "  -6.49%  "
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("D21").Value = "6.31"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("D22").Value = "65.36"
$ws.Range("E22").Value = "  -1.54%  "
$ws.Range("D23").Value = "238.73"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("E24").Value = "  -6.83%  "
$ws.Range("E25").Value = "  -8.44%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("D28").Value = "37.68"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("E29").Value = "  -4.72%  "
$ws.Range("D30").Value = "6.01"
$ws.Range("E30").Value = "  -7.56%  "
$ws.Range("D31").Value = "20.07"
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("D32").Value = "155.51"
$ws.Range("E32").Value = "  -4.12%  "
$ws.Range("D33").Value = "0.0837"
$ws.Range("E33").Value = "  -5.42%  "
$ws.Range("E34").Value = "  +8.54%  "
$ws.Range("E35").Value = "  -3.94%  "
$ws.Range("E36").Value = "  -4.95%  "
$ws.Range("D37").Value = "1.90"
$ws.Range("E37").Value = "  -6.26%  "
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("D39").Value = "15.19"
$ws.Range("E39").Value = "  -2.76%  "
$ws.Range("E40").Value = "  -10.67%  "
$ws.Range("D41").Value = "3.93"
$ws.Range("E41").Value = "  -11.62%  "
$ws.Range("D42").Value = "0.0308"
$ws.Range("E42").Value = "  -6.49%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "1.715.13"
$ws.Range("E44").Value = "  -3.30%  "
$ws.Range("D45").Value = "83.26"
$ws.Range("E45").Value = "  -4.02%  "
$ws.Range("E46").Value = "  -6.68%  "
$ws.Range("E47").Value = "  -4.99%  "
$ws.Range("E48").Value = "  -2.56%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "15.00"
$ws.Range("E49").Value = "  +3.89%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "56.51"
$ws.Range("E50").Value = "  -6.99%  "
$ws.Range("D51").Value = "1.60"
$ws.Range("E51").Value = "  -6.04%  "
